$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 154, pushing the existing data
# (previously rows 154-228) down to rows 156-230.
$ws.Rows("154:155").Insert()

# Populate new row 154 with the new weekly price record.
$ws.Range("A154").Value = 3
$ws.Range("B154").Value = "Femacal de La Calera"
$ws.Range("C154").Value = "Coquimbo"
$ws.Range("D154").Value = 44510
$ws.Range("E154").Value = 5
$ws.Range("F154").Value = 100114013
$ws.Range("G154").Value = "Zanahoria"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 180
$ws.Range("K154").Value = 7000
$ws.Range("L154").Value = 7000
$ws.Range("M154").Value = 7000
$ws.Range("N154").Value = "`$/saco 20 kilos"
$ws.Range("O154").Value = "Provincia de Quillota"
$ws.Range("P154").Value = 350
$ws.Range("Q154").Value = 20
$ws.Range("R154").Value = "Hortaliza"

# Populate new row 155 with the second new weekly price record.
$ws.Range("A155").Value = 3
$ws.Range("B155").Value = "Femacal de La Calera"
$ws.Range("C155").Value = "Coquimbo"
$ws.Range("D155").Value = 44510
$ws.Range("E155").Value = 5
$ws.Range("F155").Value = 100114013
$ws.Range("G155").Value = "Zanahoria"
$ws.Range("H155").Value = "Sin especificar"
$ws.Range("I155").Value = "Segunda"
$ws.Range("J155").Value = 160
$ws.Range("K155").Value = 5000
$ws.Range("L155").Value = 5000
$ws.Range("M155").Value = 5000
$ws.Range("N155").Value = "`$/saco 20 kilos"
$ws.Range("O155").Value = "Provincia de Quillota"
$ws.Range("P155").Value = 250
$ws.Range("Q155").Value = 20
$ws.Range("R155").Value = "Hortaliza"
